$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2382100373506546
$ws.Range("B1").Value = 0.2247202843427658
$ws.Range("C1").Value = 0.2253018766641617
$ws.Range("D1").Value = 0.2872776985168457
$ws.Range("E1").Value = 0.4499534964561462
